# Update benchmark: 2025-12-02 06:41:41 UTC
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BENCHMARK")

# Rows where column C value is cleared and moved (same or updated text) into column D.
# C is always cleared; D receives the (possibly updated) value.
$updates = @{
    3  = @{ D = "30,46 TL - 60,94 TL - 609,43 TL" }
    4  = @{ D = "30,46 TL - 60,94 TL - 609,43 TL" }
    5  = @{ D = "30,46 TL - 60,94 TL - 609,43 TL" }
    6  = @{ D = "6,09 TL - 12,19 TL - 152,35 TL" }
    8  = @{ D = "15,23 TL - 30,47 TL - 304,71 TL" }
    9  = @{ D = "15,23 TL - 30,47 TL - 304,71 TL" }
    10 = @{ D = "15,23 TL - 30,47 TL - 304,71 TL" }
    11 = @{ D = "3,04 TL - 6,09 TL - 76,17 TL" }
    12 = @{ D = "Diğer: 700 TL–4.000 TL" }
    13 = @{ D = "Hesaba: Asgari 1 TL | Azami 909,5 TL" }
    14 = @{ D = "2.300 TL - 9.500 TL" }
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 3).Value = ""
    $ws.Cells.Item($row, 4).Value = $updates[$row].D
}

# New values added to column K for rows 24 and 25 (previously empty).
$ws.Cells.Item(24, 11).Value = "371,72 TL"
$ws.Cells.Item(25, 11).Value = "312 TL"
